# Add the new "2023-24" row to the positive destinations table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row right below the header (row 2), pushing the existing
# data rows (and their formatting) down by one.
$ws.Rows(3).Insert()

# Grow the table to include the freshly inserted row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A2:H17"))

# Populate the new row with the 2023-24 figures.
$ws.Range("A3").Value2 = "2023-24"
$ws.Range("B3").Value2 = 92.78447417892292
$ws.Range("C3").Value2 = 92.559508957532813
$ws.Range("D3").Value2 = 94.098169068952089
$ws.Range("E3").Value2 = 92.578125
$ws.Range("F3").Value2 = 94.119854873013892
$ws.Range("G3").Value2 = 94.83491337038248
$ws.Range("H3").Value2 = 93.12228051641091

# Shrink the whole table's font from 12pt to 10pt so the extra row fits.
$ws.Range("A2:H17").Font.Size = 10

# New row: bold "Year" style with a left border on the year cell ...
$ws.Range("A3").Font.Bold = $true
$ws.Range("A3").VerticalAlignment = -4108
$ws.Range("A3").Borders.Item(7).LineStyle = 1
$ws.Range("A3").Borders.Item(7).Weight = 2

# ... and right-aligned whole numbers for the data cells.
$ws.Range("B3:H3").NumberFormat = "#,##0"
$ws.Range("B3:H3").HorizontalAlignment = -4152
$ws.Range("B3:H3").VerticalAlignment = -4108
$ws.Range("H3").Borders.Item(10).LineStyle = 1
$ws.Range("H3").Borders.Item(10).Weight = 2

# Highlight the new row in red if any destination rate drops below 5%.
$ws.Range("B3:H3").FormatConditions.Add(1, 6, "0.05") | Out-Null
$fc = $ws.Range("B3:H3").FormatConditions.Item($ws.Range("B3:H3").FormatConditions.Count)
$fc.Font.Color = 255

# Keep the sheet selection pointed at the (now larger) table, matching
# where Excel leaves the cursor after a table edit like this.
$ws.Range("A2:H17").Select() | Out-Null
